$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The diagonal cells (A2, B3, C4, D5) were numeric "1" values; they need to
# become the literal text "0.5" (t="s" pointing at a new shared string),
# while keeping their existing cell style untouched.
$cells = @("A2", "B3", "C4", "D5")
foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    # Force text entry (otherwise "0.5" is auto-detected as a number), then
    # restore the original "General" number format so the cell keeps using
    # its original style instead of a brand-new text-formatted style.
    $rng.NumberFormat = "@"
    $rng.Value = "0.5"
    $rng.NumberFormat = "General"
}

# Default column width changed slightly (12.703125 -> 12.68359375)
$ws.StandardWidth = 12.68359375

# Move the active cell selection to D6
$ws.Range("D6").Select()
